$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the NroSiniestro (column E) test values on the "Hoja1" sheet to the
# corrected claim numbers. A leading apostrophe forces Excel to store each
# numeric-looking, space-padded value as literal text (matching the shared
# string table in the target workbook) instead of re-interpreting it as a
# number, while leaving each cell's existing style/number format untouched.
# The edit order below reproduces the shared-string insertion order of the
# original authoring session.
$ws.Range("E5").Value = "'1120170200967   "
$ws.Range("E6").Value = "'1220194200694    "
$ws.Range("E7").Value = "'1120194100448   "
$ws.Range("E4").Value = "'1120194100448     "
$ws.Range("E3").Value = "'1220194200694 "
$ws.Range("E8").Value = "'1220194200694    "

# Reflect the cursor position saved with the workbook.
$ws.Range("G5").Select()
